$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2617.3333
$ws.Range("I62").Value = 2467.5
$ws.Range("J62").Value = 2767.1667
$ws.Range("K62").Value = 2467.5
$ws.Range("L62").Value = 2767.1667
$ws.Range("M62").Value = -1843.5
$ws.Range("N62").Value = -4015.1667

$ws.Range("H65").Value = 2617.3333
$ws.Range("I65").Value = 2467.5
$ws.Range("J65").Value = 2767.1667
$ws.Range("K65").Value = 12337.5
$ws.Range("L65").Value = 13835.8335
$ws.Range("M65").Value = -9217.5
$ws.Range("N65").Value = -20075.8335

$ws.Range("H70").Value = 1766.6666
$ws.Range("I70").Value = 1300
$ws.Range("K70").Value = 3900
$ws.Range("M70").Value = -3630

$ws.Range("H73").Value = 1766.6666
$ws.Range("I73").Value = 1300
$ws.Range("K73").Value = 3900
$ws.Range("M73").Value = -2964

$ws.Range("H113").Value = 2263.889
$ws.Range("I113").Value = 1750.5
$ws.Range("J113").Value = 2905.625
$ws.Range("K113").Value = 1750.5
$ws.Range("L113").Value = 2905.625
$ws.Range("M113").Value = 1503.5
$ws.Range("N113").Value = -9413.625

$ws.Range("H137").Value = 796.05
$ws.Range("I137").Value = 693.86365
$ws.Range("J137").Value = 920.94446
$ws.Range("K137").Value = 2081.59095
$ws.Range("L137").Value = 2762.83338
$ws.Range("M137").Value = 468.4090500000002
$ws.Range("N137").Value = -7862.83338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4414
$ws.Range("I63").Value = 2686.2856
$ws.Range("J63").Value = 5925.75
$ws.Range("K63").Value = 2686.2856
$ws.Range("L63").Value = 5925.75
$ws.Range("M63").Value = -2000.2856
$ws.Range("N63").Value = -7297.75

$ws.Range("H66").Value = 4414
$ws.Range("I66").Value = 2686.2856
$ws.Range("J66").Value = 5925.75
$ws.Range("K66").Value = 13431.428
$ws.Range("L66").Value = 29628.75
$ws.Range("M66").Value = -9999.428
$ws.Range("N66").Value = -36492.75

$ws.Range("H132").Value = 1656.8334
$ws.Range("I132").Value = 1149.1765
$ws.Range("J132").Value = 2889.7144
$ws.Range("K132").Value = 3447.5295
$ws.Range("L132").Value = 8669.143199999999
$ws.Range("M132").Value = -917.5295000000001
$ws.Range("N132").Value = -13729.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8066416
$ws.Range("I31").Value = 10205526
$ws.Range("J31").Value = 3618.4614
$ws.Range("K31").Value = 10205526
$ws.Range("L31").Value = 3618.4614
$ws.Range("M31").Value = -10205231
$ws.Range("N31").Value = -4208.4614

$ws.Range("H34").Value = 8066416
$ws.Range("I34").Value = 10205526
$ws.Range("J34").Value = 3618.4614
$ws.Range("K34").Value = 10205526
$ws.Range("L34").Value = 3618.4614
$ws.Range("M34").Value = -10205324
$ws.Range("N34").Value = -4022.4614

$ws.Range("H58").Value = 869.44446
$ws.Range("I58").Value = 689.56665
$ws.Range("J58").Value = 1229.2
$ws.Range("K58").Value = 689.56665
$ws.Range("L58").Value = 1229.2
$ws.Range("M58").Value = -486.56665
$ws.Range("N58").Value = -1635.2

$ws.Range("H136").Value = 869.44446
$ws.Range("I136").Value = 689.56665
$ws.Range("J136").Value = 1229.2
$ws.Range("K136").Value = 2068.69995
$ws.Range("L136").Value = 3687.6
$ws.Range("M136").Value = 481.3000499999998
$ws.Range("N136").Value = -8787.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 140.5625
$ws.Range("I33").Value = 60.42857
$ws.Range("J33").Value = 356.30768
$ws.Range("K33").Value = 362.57142
$ws.Range("L33").Value = 2137.84608
$ws.Range("M33").Value = -79.57141999999999
$ws.Range("N33").Value = -2703.84608

$ws.Range("H102").Value = 5307.6924
$ws.Range("J102").Value = 5307.6924
$ws.Range("L102").Value = 15923.0772
$ws.Range("N102").Value = -20791.0772

$ws.Range("H132").Value = 2555.32
$ws.Range("I132").Value = 1571.4286
$ws.Range("J132").Value = 2937.9443
$ws.Range("K132").Value = 14142.8574
$ws.Range("L132").Value = 26441.4987
$ws.Range("M132").Value = -11612.8574
$ws.Range("N132").Value = -31501.4987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 39185.594
$ws.Range("I132").Value = 49749.953
$ws.Range("J132").Value = 2210.3333
$ws.Range("K132").Value = 149249.859
$ws.Range("L132").Value = 6630.999899999999
$ws.Range("M132").Value = -146719.859
$ws.Range("N132").Value = -11690.9999

$ws.Range("H133").Value = 51282.855
$ws.Range("J133").Value = 51282.855
$ws.Range("L133").Value = 51282.855
$ws.Range("N133").Value = -61402.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 643.4909
$ws.Range("I22").Value = 506.79544
$ws.Range("J22").Value = 1190.2727
$ws.Range("K22").Value = 506.79544
$ws.Range("L22").Value = 1190.2727
$ws.Range("M22").Value = -211.79544
$ws.Range("N22").Value = -1780.2727

$ws.Range("H27").Value = 643.4909
$ws.Range("I27").Value = 506.79544
$ws.Range("J27").Value = 1190.2727
$ws.Range("K27").Value = 506.79544
$ws.Range("L27").Value = 1190.2727
$ws.Range("M27").Value = -399.79544
$ws.Range("N27").Value = -1404.2727

$ws.Range("H40").Value = 1791.4667
$ws.Range("I40").Value = 1633.8948
$ws.Range("J40").Value = 2063.6365
$ws.Range("K40").Value = 1633.8948
$ws.Range("L40").Value = 2063.6365
$ws.Range("M40").Value = -1497.8948
$ws.Range("N40").Value = -2335.6365

$ws.Range("H46").Value = 6368.8
$ws.Range("J46").Value = 10817.4
$ws.Range("L46").Value = 10817.4
$ws.Range("N46").Value = -11193.4

$ws.Range("H132").Value = 3077.5
$ws.Range("I132").Value = 3643.8333
$ws.Range("J132").Value = 1863.9286
$ws.Range("K132").Value = 10931.4999
$ws.Range("L132").Value = 5591.7858
$ws.Range("M132").Value = -8401.499899999999
$ws.Range("N132").Value = -10651.7858

$ws.Range("H136").Value = 2506
$ws.Range("I136").Value = 1594.3077
$ws.Range("J136").Value = 3987.5
$ws.Range("K136").Value = 4782.9231
$ws.Range("L136").Value = 11962.5
$ws.Range("M136").Value = -2232.9231
$ws.Range("N136").Value = -17062.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 38462930
$ws.Range("I126").Value = 76924140
$ws.Range("J126").Value = 1713.9231
$ws.Range("K126").Value = 230772420
$ws.Range("L126").Value = 5141.7693
$ws.Range("M126").Value = -230769950
$ws.Range("N126").Value = -10081.7693

$ws.Range("H136").Value = 3133.45
$ws.Range("I136").Value = 4306.269
$ws.Range("J136").Value = 955.3570999999999
$ws.Range("K136").Value = 12918.807
$ws.Range("L136").Value = 2866.0713
$ws.Range("M136").Value = -10368.807
$ws.Range("N136").Value = -7966.0713
